# Refresh the "last updated" timestamp column (F) on each data sheet.
# The JSON database finished building, so every status row gets a new
# last-updated stamp reflecting the latest refresh run.

$wb = $excel.ActiveWorkbook

# Sheet 1: "שרתים" (Servers) - rows 2-7
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2:F7").Value = "21/10/2021, 11:29:41"

# Sheet 2: "תא קשר" (Contact cell) - rows 2-22 (row 23 keeps its older stamp)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2:F22").Value = "21/10/2021, 11:29:24"

# Sheet 3: "קרונות" (Carriages) - rows 2-32
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2:F32").Value = "21/10/2021, 11:29:30"

# Sheet 4: "לבדוק בנוסף" (To check additionally) - only rows 11, 12, 17
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F11:F12").Value = "21/10/2021, 11:29:13"
$ws4.Range("F17").Value = "21/10/2021, 11:29:13"
